$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 50555.1
$ws.Range("I28").Value = 59266.65
$ws.Range("K28").Value = 59266.65
$ws.Range("M28").Value = -58781.65
$ws.Range("H29").Value = 5541.75
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 5541.75
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 16625.25
$ws.Range("M29").Value = $null
$ws.Range("N29").Value = -17187.25
$ws.Range("H64").Value = 7186.25
$ws.Range("I64").Value = 4247.5
$ws.Range("J64").Value = 10125
$ws.Range("K64").Value = 4247.5
$ws.Range("L64").Value = 10125
$ws.Range("M64").Value = -3999.5
$ws.Range("N64").Value = -10621
$ws.Range("H67").Value = 7186.25
$ws.Range("I67").Value = 4247.5
$ws.Range("J67").Value = 10125
$ws.Range("K67").Value = 4247.5
$ws.Range("L67").Value = 10125
$ws.Range("M67").Value = -3389.5
$ws.Range("N67").Value = -11841
$ws.Range("H132").Value = 4011.0732
$ws.Range("I132").Value = 3986.2
$ws.Range("K132").Value = 11958.6
$ws.Range("M132").Value = -9428.599999999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 970.79486
$ws.Range("I32").Value = 969.6667
$ws.Range("K32").Value = 969.6667
$ws.Range("M32").Value = -682.6667
$ws.Range("H102").Value = 2261.6086
$ws.Range("I102").Value = 2261.6086
$ws.Range("K102").Value = 2261.6086
$ws.Range("M102").Value = -639.6086

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 9260935
$ws.Range("I105").Value = 1830.7858
$ws.Range("K105").Value = 1830.7858
$ws.Range("M105").Value = -83.78580000000011

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1054.5264
$ws.Range("I16").Value = 977.93335
$ws.Range("K16").Value = 977.93335
$ws.Range("M16").Value = -690.93335
$ws.Range("H31").Value = 792239.2
$ws.Range("I31").Value = 1541969.1
$ws.Range("J31").Value = 42509.184
$ws.Range("K31").Value = 1541969.1
$ws.Range("L31").Value = 42509.184
$ws.Range("M31").Value = -1541674.1
$ws.Range("N31").Value = -43099.184
$ws.Range("H34").Value = 792239.2
$ws.Range("I34").Value = 1541969.1
$ws.Range("J34").Value = 42509.184
$ws.Range("K34").Value = 1541969.1
$ws.Range("L34").Value = 42509.184
$ws.Range("M34").Value = -1541767.1
$ws.Range("N34").Value = -42913.184
$ws.Range("H99").Value = 5849.75
$ws.Range("I99").Value = 3533
$ws.Range("K99").Value = 3533
$ws.Range("M99").Value = -2035
$ws.Range("H113").Value = 1054.5264
$ws.Range("I113").Value = 977.93335
$ws.Range("K113").Value = 977.93335
$ws.Range("M113").Value = 1192.06665
$ws.Range("H122").Value = 2319.238
$ws.Range("I122").Value = 1850
$ws.Range("J122").Value = 3257.7144
$ws.Range("K122").Value = 5550
$ws.Range("L122").Value = 9773.143199999999
$ws.Range("M122").Value = -3100
$ws.Range("N122").Value = -14673.1432
$ws.Range("H126").Value = 5849.75
$ws.Range("I126").Value = 3533
$ws.Range("K126").Value = 10599
$ws.Range("M126").Value = -8129

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 4174.033
$ws.Range("I14").Value = 4174.033
$ws.Range("K14").Value = 12522.099
$ws.Range("M14").Value = -12349.099
$ws.Range("H121").Value = 909634
$ws.Range("I121").Value = 445.8
$ws.Range("J121").Value = 1667290.9
$ws.Range("K121").Value = 1337.4
$ws.Range("L121").Value = 5001872.699999999
$ws.Range("M121").Value = -27.40000000000009
$ws.Range("N121").Value = -5004492.699999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 1324.2858
$ws.Range("I17").Value = 741
$ws.Range("J17").Value = 2102
$ws.Range("K17").Value = 741
$ws.Range("L17").Value = 2102
$ws.Range("M17").Value = -573
$ws.Range("N17").Value = -2438
$ws.Range("H24").Value = 50000
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 50000
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 50000
$ws.Range("M24").Value = $null
$ws.Range("N24").Value = -50346
$ws.Range("H97").Value = 602.4666999999999
$ws.Range("I97").Value = 669.8182
$ws.Range("J97").Value = 417.25
$ws.Range("K97").Value = 669.8182
$ws.Range("L97").Value = 417.25
$ws.Range("M97").Value = -173.8182
$ws.Range("N97").Value = -1409.25
$ws.Range("H113").Value = 59444784
$ws.Range("I113").Value = 2091468.6
$ws.Range("J113").Value = 83342000
$ws.Range("K113").Value = 2091468.6
$ws.Range("L113").Value = 83342000
$ws.Range("M113").Value = -2089298.6
$ws.Range("N113").Value = -83346340
$ws.Range("H132").Value = 45380
$ws.Range("I132").Value = 15968.125
$ws.Range("J132").Value = 146220.72
$ws.Range("K132").Value = 47904.375
$ws.Range("L132").Value = 438662.16
$ws.Range("M132").Value = -45374.375
$ws.Range("N132").Value = -443722.16

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 444773.44
$ws.Range("I7").Value = 719286.3
$ws.Range("J7").Value = 17753.445
$ws.Range("K7").Value = 719286.3
$ws.Range("L7").Value = 17753.445
$ws.Range("M7").Value = -719174.3
$ws.Range("N7").Value = -17977.445
$ws.Range("H13").Value = 9666.666999999999
$ws.Range("I13").Value = 9000
$ws.Range("J13").Value = 10000
$ws.Range("K13").Value = 9000
$ws.Range("L13").Value = 10000
$ws.Range("M13").Value = -8860
$ws.Range("N13").Value = -10280
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = $null
$ws.Range("N23").Value = $null
$ws.Range("H26").Value = 18000
$ws.Range("J26").Value = 18000
$ws.Range("L26").Value = 18000
$ws.Range("N26").Value = -18590
$ws.Range("H31").Value = 2772.8572
$ws.Range("I31").Value = 957.5
$ws.Range("J31").Value = 3499
$ws.Range("K31").Value = 957.5
$ws.Range("L31").Value = 3499
$ws.Range("M31").Value = -709.5
$ws.Range("N31").Value = -3995
$ws.Range("H38").Value = 29976.666
$ws.Range("I38").Value = 24965
$ws.Range("K38").Value = 24965
$ws.Range("M38").Value = -24555
$ws.Range("H46").Value = 2693.853
$ws.Range("I46").Value = 2510.158
$ws.Range("J46").Value = 2926.5334
$ws.Range("K46").Value = 2510.158
$ws.Range("L46").Value = 2926.5334
$ws.Range("M46").Value = -2322.158
$ws.Range("N46").Value = -3302.5334
$ws.Range("H61").Value = 3573.0908
$ws.Range("I61").Value = 2666.5833
$ws.Range("K61").Value = 2666.5833
$ws.Range("M61").Value = -2464.5833
$ws.Range("H113").Value = 3573.0908
$ws.Range("I113").Value = 2666.5833
$ws.Range("K113").Value = 2666.5833
$ws.Range("M113").Value = -496.5832999999998
$ws.Range("H122").Value = 583777.7
$ws.Range("I122").Value = 3419.6667
$ws.Range("K122").Value = 10259.0001
$ws.Range("M122").Value = -7809.000100000001
$ws.Range("H126").Value = 444773.44
$ws.Range("I126").Value = 719286.3
$ws.Range("J126").Value = 17753.445
$ws.Range("K126").Value = 2157858.9
$ws.Range("L126").Value = 53260.335
$ws.Range("M126").Value = -2155388.9
$ws.Range("N126").Value = -58200.335

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1802.9412
$ws.Range("I113").Value = 1876.909
$ws.Range("K113").Value = 5630.727000000001
$ws.Range("M113").Value = -3460.727000000001
